# This script completes the partially-populated sweep groups (B = 8, 9, 10)
# in the simulation results sheet by inserting the missing Laser-Power (A)
# rows in their correct sorted position, each carrying the constant
# per-scenario columns D:H (0, 250, 112, 0.16, 0.158).
#
# Rows are inserted from the bottom of the range upward so that earlier
# (lower) insertion points are not invalidated by later inserts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-DataRow($BeforeRow, $A, $B, $C) {
    $ws.Rows.Item($BeforeRow).Insert()

    $ws.Cells.Item($BeforeRow, 1).Value = $A
    $ws.Cells.Item($BeforeRow, 2).Value = $B
    $ws.Cells.Item($BeforeRow, 3).Value = $C
    $ws.Cells.Item($BeforeRow, 4).Value = 0
    $ws.Cells.Item($BeforeRow, 5).Value = 250
    $ws.Cells.Item($BeforeRow, 6).Value = 112
    $ws.Cells.Item($BeforeRow, 7).Value = 0.16
    $ws.Cells.Item($BeforeRow, 8).Value = 0.158
}

# --- Insert from the highest original row position down to the lowest ---

# Before original row 113 (old B=10 group): A = -6
Insert-DataRow 113 -6 10 5.9466

# Before original row 111 (old B=9/B=10 boundary): A = 0, 1, 2 (B=9) then A = -9 (B=10)
Insert-DataRow 111 0  9  4.6252
Insert-DataRow 112 1  9  3.544
Insert-DataRow 113 2  9  1.4599
Insert-DataRow 114 -9 10 4.3695

# Before original row 108 (old B=9 group): A = -4
Insert-DataRow 108 -4 9 6.6913

# Before original row 107 (old B=9 group): A = -8, -7, -6
Insert-DataRow 107 -8 9 4.9847
Insert-DataRow 108 -7 9 5.4923
Insert-DataRow 109 -6 9 5.964

# Before original row 105 (old B=8 group): A = 0
Insert-DataRow 105 0 8 4.2935
